$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 739.58826
$ws.Range("I28").Value = 754.5625
$ws.Range("J28").Value = 500
$ws.Range("K28").Value = 754.5625
$ws.Range("L28").Value = 500
$ws.Range("M28").Value = -269.5625
$ws.Range("N28").Value = -1470

# Row 111
$ws.Range("H111").Value = 1365
$ws.Range("I111").Value = 1345
$ws.Range("J111").Value = 1375
$ws.Range("K111").Value = 4035
$ws.Range("L111").Value = 4125
$ws.Range("M111").Value = -968
$ws.Range("N111").Value = -10259

# Row 113
$ws.Range("H113").Value = 2500
$ws.Range("I113").Value = 2000
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 2000
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 1254
$ws.Range("N113").Value = -9508

# Row 115
$ws.Range("H115").Value = 1139.5
$ws.Range("I115").Value = 1139.5
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 3418.5
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = -1851.5

# Row 138
$ws.Range("H138").Value = 2468.5
$ws.Range("I138").Value = 1335.375
$ws.Range("J138").Value = 3035.0625
$ws.Range("K138").Value = 4006.125
$ws.Range("L138").Value = 9105.1875
$ws.Range("M138").Value = 1133.875
$ws.Range("N138").Value = -19385.1875

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 794.1667
$ws.Range("I2").Value = 794.1667
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 794.1667
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -681.1667
$ws.Range("N2").ClearContents()

# Row 94
$ws.Range("H94").Value = 54999
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 54999
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 54999
$ws.Range("N94").Value = -56801

# Row 110
$ws.Range("H110").Value = 15625613
$ws.Range("I110").Value = 570.1111
$ws.Range("J110").Value = 35714950
$ws.Range("K110").Value = 570.1111
$ws.Range("L110").Value = 35714950
$ws.Range("M110").Value = 1474.8889
$ws.Range("N110").Value = -35719040

# Row 116
$ws.Range("H116").Value = 794.1667
$ws.Range("I116").Value = 794.1667
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 794.1667
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1499.8333
$ws.Range("N116").ClearContents()

# Row 117
$ws.Range("H117").Value = 76000
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 76000
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 76000
$ws.Range("N117").Value = -85178

# Row 120
$ws.Range("H120").Value = 0
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()

# Row 121
$ws.Range("H121").Value = 0
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 794.1667
$ws.Range("I3").Value = 794.1667
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 794.1667
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -680.1667
$ws.Range("N3").ClearContents()

# Row 46
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("N46").ClearContents()

# Row 86
$ws.Range("H86").Value = 1200.75
$ws.Range("I86").Value = 1184.3334
$ws.Range("J86").Value = 1250
$ws.Range("K86").Value = 1184.3334
$ws.Range("L86").Value = 1250
$ws.Range("M86").Value = -61.33339999999998
$ws.Range("N86").Value = -3496

# Row 89
$ws.Range("H89").Value = 1200.75
$ws.Range("I89").Value = 1184.3334
$ws.Range("J89").Value = 1250
$ws.Range("K89").Value = 5921.666999999999
$ws.Range("L89").Value = 6250
$ws.Range("M89").Value = -305.6669999999995
$ws.Range("N89").Value = -17482

# Row 100
$ws.Range("H100").Value = 13215
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 13215
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 13215
$ws.Range("N100").Value = -15379

# Row 107
$ws.Range("H107").Value = 400000
$ws.Range("I107").Value = 400000
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 400000
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -398080

# Row 134
$ws.Range("H134").Value = 1019.0833
$ws.Range("I134").Value = 1019.0833
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 3057.2499
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -522.2498999999998

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 767.7368
$ws.Range("I22").Value = 760.3889
$ws.Range("J22").Value = 900
$ws.Range("K22").Value = 760.3889
$ws.Range("L22").Value = 900
$ws.Range("M22").Value = -410.3889
$ws.Range("N22").Value = -1600

# Row 107
$ws.Range("H107").Value = 866.3333
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 866.3333
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 866.3333
$ws.Range("N107").Value = -4706.3333

$ws = $wb.Worksheets.Item("CUL")
# Row 39
$ws.Range("H39").Value = 14222.223
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 14222.223
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 42666.669
$ws.Range("N39").Value = -43254.669

# Row 92
$ws.Range("H92").Value = 1097.6666
$ws.Range("I92").Value = 700
$ws.Range("J92").Value = 1296.5
$ws.Range("K92").Value = 2100
$ws.Range("L92").Value = 3889.5
$ws.Range("M92").Value = -852
$ws.Range("N92").Value = -6385.5

# Row 114
$ws.Range("H114").Value = 1282
$ws.Range("I114").Value = 1376
$ws.Range("J114").Value = 1000
$ws.Range("K114").Value = 4128
$ws.Range("L114").Value = 3000
$ws.Range("M114").Value = -874
$ws.Range("N114").Value = -9508

# Row 129
$ws.Range("H129").Value = 1398.2
$ws.Range("I129").Value = 1029
$ws.Range("J129").Value = 1644.3334
$ws.Range("K129").Value = 3087
$ws.Range("L129").Value = 4933.0002
$ws.Range("M129").Value = 1913
$ws.Range("N129").Value = -14933.0002

$ws = $wb.Worksheets.Item("GSM")
# Row 93
$ws.Range("H93").Value = 31666.666
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 31666.666
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 31666.666
$ws.Range("N93").Value = -35410.666

# Row 113
$ws.Range("H113").Value = 2874.6667
$ws.Range("I113").Value = 1812.25
$ws.Range("J113").Value = 4999.5
$ws.Range("K113").Value = 1812.25
$ws.Range("L113").Value = 4999.5
$ws.Range("M113").Value = 357.75
$ws.Range("N113").Value = -9339.5

# Row 132
$ws.Range("H132").Value = 1326.3334
$ws.Range("I132").Value = 1326.3334
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3979.0002
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1449.0002
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 8100.3335
$ws.Range("I16").Value = 4333.3335
$ws.Range("J16").Value = 11867.333
$ws.Range("K16").Value = 4333.3335
$ws.Range("L16").Value = 11867.333
$ws.Range("M16").Value = -4163.3335
$ws.Range("N16").Value = -12207.333

# Row 82
$ws.Range("H82").Value = 1460
$ws.Range("I82").Value = 1460
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 1460
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -1099
$ws.Range("N82").ClearContents()

# Row 85
$ws.Range("H85").Value = 1460
$ws.Range("I85").Value = 1460
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 1460
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -212
$ws.Range("N85").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 39
$ws.Range("H39").Value = 10000
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 10000
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 10000
$ws.Range("N39").Value = -10826

Write-Host "Applied Golem_Profits updates"
